# Auto-committed on 2022/04/28 週四
# Adds a new GenTable entry "CdConvertCode" (代碼轉換檔) into the
# "L6-共同作業" section of the table (alphabetically between CdCl and CdEmp),
# and corrects the "AML定審資料" last-modified timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right before the existing "CdEmp" row (row 160), which
# pushes every following row down by one and keeps all formulas/styles intact.
$ws.Rows(160).Insert()

# Fill in the new row for the CdConvertCode table.
$ws.Range("A160").Value = "L6-共同作業"
$ws.Range("B160").Value = "CdConvertCode"
$ws.Range("C160").Value = "代碼轉換檔"
$ws.Range("D160").Formula = "=HYPERLINK(""[\\192.168.10.16\St1Share(NAS)\SKL\DB\GenTables\L6-共同作業\CdConvertCode.xlsx]DBD!A1"", ""連結"")"
$ws.Range("E160").Value = "2022年04月28日 16:39:43"

# Correct the last-modified timestamp recorded for the "TxAmlCredit"
# (AML定審資料) table row, which shifted from row 323 to row 324 because of
# the row insertion above.
$ws.Range("E324").Value = "2022年04月28日 14:17:55"
